# Generate Report for Handback
#
# The handback status report records, per locale sheet, the timestamps at
# which each handed-off file's round-trip (handoff -> handback) completed.
# A new report run completed processing for the "b5a67361-..." file (the
# table's row 2) in both the "zh-cn" and "de-de" locale sheets, so its
# "Correspond Handoff Datetime" (column E) and "Correspond Handback
# DateTime" (column H) are refreshed to the new run's timestamps. The
# "f84312db-..." file (row 3) was not part of this run, so its timestamps
# are left untouched.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 12:52:45"
$wsZhCn.Range("H2").Value = "2016-03-23 12:53:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 12:52:49"
$wsDeDe.Range("H2").Value = "2016-03-23 12:53:22"
